# "Generate Report for Handoff"
#
# The two tracked files (83f153c1-... and c88d606e-...) swap status rows:
#   - c88d606e-...  was already "Handed back: in sync with en-US" and stays
#     that way (it now occupies row 2 on every sheet).
#   - 83f153c1-...  moves to row 3 and is updated to "Ready for handoff"
#     with fresh handoff/target timestamps (it was previously handed back
#     too, in row 2).
#
# Underlying hyperlink relationships (r:id -> external github URL) are left
# untouched; only the cell values and each hyperlink's displayed text are
# updated, exactly mirroring what Excel does when you edit .Value /
# .TextToDisplay without re-inserting the hyperlink.

$wb = $excel.ActiveWorkbook

function Set-CellAndMaybeHyperlink {
    param(
        $ws,
        [string]$addr,
        [string]$value
    )
    $ws.Range($addr).Value = $value
    foreach ($hl in $ws.Hyperlinks) {
        $hlAddr = $hl.Range.Address()
        if ($hlAddr -eq ('$' + $addr.Substring(0,1) + '$' + $addr.Substring(1))) {
            $hl.TextToDisplay = $value
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-CellAndMaybeHyperlink $ws "A2" "c88d606e-1676-4d16-a30c-f85b4acc0204.md"
Set-CellAndMaybeHyperlink $ws "A3" "83f153c1-74e8-465e-9c49-b796e5acc545.md"

$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-22 04:46:11"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-CellAndMaybeHyperlink $ws "A2" "c88d606e-1676-4d16-a30c-f85b4acc0204.md"
Set-CellAndMaybeHyperlink $ws "D2" "c88d606e-1676-4d16-a30c-f85b4acc0204.c612a63089d6809fbc64fc68fe4fca339d8eb1cf.zh-cn.xlf"
Set-CellAndMaybeHyperlink $ws "F2" "c88d606e-1676-4d16-a30c-f85b4acc0204.md"
Set-CellAndMaybeHyperlink $ws "G2" "c88d606e-1676-4d16-a30c-f85b4acc0204.c612a63089d6809fbc64fc68fe4fca339d8eb1cf.zh-cn.xlf"

Set-CellAndMaybeHyperlink $ws "A3" "83f153c1-74e8-465e-9c49-b796e5acc545.md"
Set-CellAndMaybeHyperlink $ws "D3" "83f153c1-74e8-465e-9c49-b796e5acc545.f1beab7333bf4a583230ef556786908b92f927a8.zh-cn.xlf"
Set-CellAndMaybeHyperlink $ws "F3" "83f153c1-74e8-465e-9c49-b796e5acc545.md"
Set-CellAndMaybeHyperlink $ws "G3" "83f153c1-74e8-465e-9c49-b796e5acc545.f1beab7333bf4a583230ef556786908b92f927a8.zh-cn.xlf"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-22 04:46:07"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-CellAndMaybeHyperlink $ws "A2" "c88d606e-1676-4d16-a30c-f85b4acc0204.md"
Set-CellAndMaybeHyperlink $ws "D2" "c88d606e-1676-4d16-a30c-f85b4acc0204.c612a63089d6809fbc64fc68fe4fca339d8eb1cf.de-de.xlf"
Set-CellAndMaybeHyperlink $ws "F2" "c88d606e-1676-4d16-a30c-f85b4acc0204.md"
Set-CellAndMaybeHyperlink $ws "G2" "c88d606e-1676-4d16-a30c-f85b4acc0204.c612a63089d6809fbc64fc68fe4fca339d8eb1cf.de-de.xlf"

Set-CellAndMaybeHyperlink $ws "A3" "83f153c1-74e8-465e-9c49-b796e5acc545.md"
Set-CellAndMaybeHyperlink $ws "D3" "83f153c1-74e8-465e-9c49-b796e5acc545.f1beab7333bf4a583230ef556786908b92f927a8.de-de.xlf"
Set-CellAndMaybeHyperlink $ws "F3" "83f153c1-74e8-465e-9c49-b796e5acc545.md"
Set-CellAndMaybeHyperlink $ws "G3" "83f153c1-74e8-465e-9c49-b796e5acc545.f1beab7333bf4a583230ef556786908b92f927a8.de-de.xlf"

$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-22 04:46:11"
